$d = $word.ActiveDocument

# 1. Remove the stray _GoBack bookmark that currently sits after "Aramis"
#    in the Question 1 paragraph (it will be re-created at the end of the
#    new Question 2 paragraph below).
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# 2. Fill in the (currently empty) Question 2 list paragraph with the new
#    question text, a line break, and the bolded answer, finishing with a
#    fresh _GoBack bookmark - mirroring exactly what Word leaves behind
#    after a user types new content at the end of the document.
$p4 = $d.Paragraphs.Item(4)
$r = $p4.Range

$xml = @'
<w:p xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
  <w:pPr>
    <w:pStyle w:val="ListParagraph"/>
    <w:numPr>
      <w:ilvl w:val="0"/>
      <w:numId w:val="3"/>
    </w:numPr>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
  </w:pPr>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>Berapakah</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>ketinggian</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>puncak</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>tertinggi</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t xml:space="preserve"> di </w:t>
  </w:r>
  <w:proofErr w:type="spellStart"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>dunia</w:t>
  </w:r>
  <w:proofErr w:type="spellEnd"/>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>, Mount Everest?</w:t>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:br/>
  </w:r>
  <w:r>
    <w:rPr>
      <w:rFonts w:ascii="Aparajita" w:hAnsi="Aparajita" w:cs="Aparajita"/>
      <w:b/>
      <w:sz w:val="24"/>
      <w:szCs w:val="24"/>
    </w:rPr>
    <w:t>8,848 Metres</w:t>
  </w:r>
  <w:bookmarkStart w:id="0" w:name="_GoBack"/>
  <w:bookmarkEnd w:id="0"/>
</w:p>
'@

$r.InsertXML($xml)
